$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Already-merged groups grow to absorb a newly-matching neighbour row:
# B5:B6 -> B5:B7 (drops the now-redundant duplicate "Three" in B7)
$ws.Range("B5:B7").Merge() | Out-Null
# E2:E3 -> E2:E4
$ws.Range("E2:E4").Merge() | Out-Null

# New grid of grouping indicator values in columns G:L, rows 1-10
$ws.Range("G1").Value = "One"
$ws.Range("H1").Value = "Two"
$ws.Range("I1").Value = "One"
$ws.Range("J1").Value = "One"
$ws.Range("K1").Value = "Two"
$ws.Range("L1").Value = "One"

$ws.Range("G2").Value = "One"
$ws.Range("H2").Value = "Two"
$ws.Range("I2").Value = "One"
$ws.Range("J2").Value = "One"
$ws.Range("K2").Value = "Two"

$ws.Range("H3").Value = "One"
$ws.Range("K3").Value = "One"
$ws.Range("L3").Value = "One"

$ws.Range("J8").Value = "Two"
$ws.Range("L8").Value = "Two"

$ws.Range("I9").Value = "Two"
$ws.Range("K9").Value = "Two"

$ws.Range("G10").Value = "One"
$ws.Range("H10").Value = "One"
$ws.Range("I10").Value = "Two"
$ws.Range("J10").Value = "Two"
$ws.Range("K10").Value = "Two"

$ws.Range("L1:L2").Merge() | Out-Null
$ws.Range("G2:G9").Merge() | Out-Null
$ws.Range("I2:I8").Merge() | Out-Null
$ws.Range("J2:J7").Merge() | Out-Null
$ws.Range("H3:H9").Merge() | Out-Null
$ws.Range("K3:K8").Merge() | Out-Null
$ws.Range("L3:L7").Merge() | Out-Null
$ws.Range("J8:J9").Merge() | Out-Null
$ws.Range("L8:L9").Merge() | Out-Null
